# Katalog guncellendi - Pzt 24.11.2025 18:05:45,99
#
# Rows 76-88 (column E) describe the "Beden" (size) options available for
# the "SISME YELEK" (puffer vest) products. They previously referenced the
# shared string "S-M-L-XL-2XL Beden secenegi mevcuttur...." ; this adds the
# 3XL size to the offering by writing the updated text into those cells.
# Excel's shared-string table will pick up the new, distinct string
# automatically (appended as a brand new entry) while leaving the old
# string intact for the other rows (58-75) that still only go up to 2XL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$yeniBedenMetni = "S-M-L-XL-2XL-3XL Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır.Belirtilen fiyatlar adet fiyatıdır."

for ($row = 76; $row -le 88; $row++) {
    $ws.Cells.Item($row, 5).Value = $yeniBedenMetni
}

# Restore/refresh the active selection to match the saved view state (E88).
$ws.Range("E88").Select()
